$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Data")

$ws.Range("B88").Formula = "=B54"
$ws.Range("C88").Formula = "=B55"
$ws.Range("B96").Formula = "=B87"

$ws.Range("B88:C88").NumberFormat = """$""#,##0"
$ws.Range("B96").NumberFormat = """$""#,##0"
